$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.516.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.034.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.88%  "

# Row 6
$ws.Range("E6").Value = "  -0.78%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.94%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0796"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "

# Row 11
$ws.Range("E11").Value = "  -1.71%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.336.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.99%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.823"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.87%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.041.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.525.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0856"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.56%  "

# Row 24
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("E25").Value = "  -1.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("E28").Value = "  -7.25%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.03%  "

# Row 30
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0668"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.19%  "

# Row 35
$ws.Range("E35").Value = "  +8.56%  "

# Row 36
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.92%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.30%  "

# Row 39
$ws.Range("E39").Value = "  -3.04%  "

# Row 40
$ws.Range("E40").Value = "  +3.87%  "

# Row 41
$ws.Range("E41").Value = "  -3.03%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0216"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.42%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.404.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.42%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.83%  "

# Row 47
$ws.Range("E47").Value = "  +1.16%  "

# Row 48
$ws.Range("E48").Value = "  +1.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.75%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.227.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.00%  "
